# Automatische test-sync: 2025-06-24 21:47:50
# Append the new "Technische storing" mail-log entry as row 37 on the
# "Logs" sheet, extend the conditional-formatting ranges to cover it,
# and bump the "IT / Technisch probleem" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A37").Value = "Technische storing"
$ws.Range("B37").Value = "mailmind.test@zohomail.eu"
$ws.Range("C37").Value = "De website werkt niet goed. Is hier iets mis mee?"
$ws.Range("D37").Value = "IT / Technisch probleem"
$ws.Range("E37").Value = "Beste klant,`nBedankt voor je bericht. Om het probleem beter te kunnen onderzoeken, hebben we meer informatie nodig. Zou je alsjeblieft kunnen aangeven welke specifieke problemen je ervaart wanneer je de website probeert te gebruiken? Bijvoorbeeld, krijg je een foutmelding te zien of lukt het niet om in te loggen? Met deze details kunnen we het probleem gericht aanpakken en een oplossing bieden.`nWe horen graag meer van je, zodat we je verder kunnen helpen.`nMet vriendelijke groet,  `n[Tekstschrijver]  `nE-mailassistent"
$ws.Range("F37").Value = "2025-06-24 21:47:18"
$ws.Range("G37").Value = "Ja"

# Writing the multi-line text into a brand-new row makes the engine pin an
# explicit (auto-estimated) row height; AutoFit clears that custom-height
# flag again so row 37 ends up with the same "no explicit height" shape as
# every other data row in the sheet.
$ws.Rows.Item(37).AutoFit()

# Extend the two conditional-formatting blocks (Categorie / Beantwoord
# columns) so they keep covering the whole data range through row 37.
$fcsCategorie = $ws.Range("D2:D36").FormatConditions
$fcsCategorie.Item(1).ModifyAppliesToRange($ws.Range("D2:D37"))

$fcsBeantwoord = $ws.Range("G2:G36").FormatConditions
$fcsBeantwoord.Item(1).ModifyAppliesToRange($ws.Range("G2:G37"))

# Update the Dashboard summary count for "IT / Technisch probleem".
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B4").Value = 5
